# Auto-generated edit script: updates market-data-derived profit columns
# (H, I, J, K, L, M, N) across multiple sheets to match refreshed values
# pulled by the scheduled market-data runner. Values are plain numeric
# literals (no formulas in this workbook), so each touched cell is set
# directly; two cells that no longer carry data are cleared entirely.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 533852.3
$ws.Range("I88").Value = 715797.9399999999
$ws.Range("J88").Value = 33502
$ws.Range("K88").Value = 715797.9399999999
$ws.Range("L88").Value = 33502
$ws.Range("M88").Value = -715391.9399999999
$ws.Range("N88").Value = -34314
$ws.Range("H91").Value = 533852.3
$ws.Range("I91").Value = 715797.9399999999
$ws.Range("J91").Value = 33502
$ws.Range("K91").Value = 715797.9399999999
$ws.Range("L91").Value = 33502
$ws.Range("M91").Value = -714393.9399999999
$ws.Range("N91").Value = -36310
$ws.Range("H108").Value = 41578.5
$ws.Range("J108").Value = 41578.5
$ws.Range("L108").Value = 41578.5
$ws.Range("N108").Value = -49258.5
$ws.Range("H112").Value = 1777.8846
$ws.Range("J112").Value = 1860.2273
$ws.Range("L112").Value = 5580.6819
$ws.Range("N112").Value = -7796.6819
$ws.Range("H130").Value = 49772
$ws.Range("J130").Value = 49772
$ws.Range("L130").Value = 49772
$ws.Range("N130").Value = -59812
$ws.Range("H138").Value = 1405.75
$ws.Range("I138").Value = 625.5
$ws.Range("J138").Value = 2186
$ws.Range("K138").Value = 1876.5
$ws.Range("L138").Value = 6558
$ws.Range("M138").Value = 3263.5
$ws.Range("N138").Value = -16838

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10772.328
$ws.Range("I32").Value = 10010.091
$ws.Range("J32").Value = 15430.444
$ws.Range("K32").Value = 10010.091
$ws.Range("L32").Value = 15430.444
$ws.Range("M32").Value = -9723.091
$ws.Range("N32").Value = -16004.444
$ws.Range("H80").Value = 47053.75
$ws.Range("J80").Value = 47053.75
$ws.Range("L80").Value = 47053.75
$ws.Range("N80").Value = -49049.75
$ws.Range("H83").Value = 47053.75
$ws.Range("J83").Value = 47053.75
$ws.Range("L83").Value = 141161.25
$ws.Range("N83").Value = -151145.25
$ws.Range("H92").Value = 16775
$ws.Range("J92").Value = 16775
$ws.Range("L92").Value = 16775
$ws.Range("N92").Value = -21767
$ws.Range("H109").Value = 46491.25
$ws.Range("J109").Value = 46491.25
$ws.Range("L109").Value = 46491.25
$ws.Range("N109").Value = -49265.25
$ws.Range("H113").Value = 46661.332
$ws.Range("J113").Value = 46661.332
$ws.Range("L113").Value = 46661.332
$ws.Range("N113").Value = -55339.332
$ws.Range("H117").Value = 41303.89
$ws.Range("J117").Value = 41303.89
$ws.Range("L117").Value = 41303.89
$ws.Range("N117").Value = -50481.89
$ws.Range("H118").Value = 49401
$ws.Range("J118").Value = 49401
$ws.Range("L118").Value = 49401
$ws.Range("N118").Value = -52715
$ws.Range("H131").Value = 48026.8
$ws.Range("J131").Value = 48026.8
$ws.Range("L131").Value = 48026.8
$ws.Range("N131").Value = -58106.8
$ws.Range("H139").Value = 49758
$ws.Range("J139").Value = 49758
$ws.Range("L139").Value = 49758
$ws.Range("N139").Value = -60038

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 31146.8
$ws.Range("J35").Value = 31146.8
$ws.Range("L35").Value = 31146.8
$ws.Range("N35").Value = -31766.8
$ws.Range("H132").Value = 45732.5
$ws.Range("J132").Value = 45732.5
$ws.Range("L132").Value = 45732.5
$ws.Range("N132").Value = -55852.5
$ws.Range("H134").Value = 2451.9355
$ws.Range("I134").Value = 1013.28
$ws.Range("K134").Value = 3039.84
$ws.Range("M134").Value = -504.8400000000001
$ws.Range("H138").Value = 195000
$ws.Range("J138").Value = 195000
$ws.Range("L138").Value = 195000
$ws.Range("N138").Value = -205280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1495.279
$ws.Range("I58").Value = 1151.1515
$ws.Range("J58").Value = 2630.9
$ws.Range("K58").Value = 1151.1515
$ws.Range("L58").Value = 2630.9
$ws.Range("M58").Value = -948.1514999999999
$ws.Range("N58").Value = -3036.9
$ws.Range("H118").Value = 47728.668
$ws.Range("J118").Value = 47728.668
$ws.Range("L118").Value = 47728.668
$ws.Range("N118").Value = -51042.668
$ws.Range("H119").Value = 42842
$ws.Range("J119").Value = 42842
$ws.Range("L119").Value = 42842
$ws.Range("N119").Value = -52518
$ws.Range("H131").Value = 36706.668
$ws.Range("J131").Value = 36706.668
$ws.Range("L131").Value = 36706.668
$ws.Range("N131").Value = -46786.668
$ws.Range("H132").Value = 35109.355
$ws.Range("I132").Value = 1589.2903
$ws.Range("K132").Value = 4767.8709
$ws.Range("M132").Value = -2237.8709
$ws.Range("H136").Value = 1495.279
$ws.Range("I136").Value = 1151.1515
$ws.Range("J136").Value = 2630.9
$ws.Range("K136").Value = 3453.4545
$ws.Range("L136").Value = 7892.700000000001
$ws.Range("M136").Value = -903.4544999999998
$ws.Range("N136").Value = -12992.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 2670.5881
$ws.Range("J49").Value = 2670.5881
$ws.Range("L49").Value = 8011.7643
$ws.Range("N49").Value = -8323.764299999999
$ws.Range("H131").Value = 2845.5254
$ws.Range("I131").Value = 9481.637000000001
$ws.Range("J131").Value = 1324.75
$ws.Range("K131").Value = 28444.911
$ws.Range("L131").Value = 3974.25
$ws.Range("M131").Value = -23404.911
$ws.Range("N131").Value = -14054.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 47659.332
$ws.Range("J110").Value = 47659.332
$ws.Range("L110").Value = 47659.332
$ws.Range("N110").Value = -55839.332
$ws.Range("H114").Value = 38559
$ws.Range("J114").Value = 38559
$ws.Range("L114").Value = 38559
$ws.Range("N114").Value = -47237
$ws.Range("H130").Value = 44956
$ws.Range("J130").Value = 44956
$ws.Range("L130").Value = 44956
$ws.Range("N130").Value = -54996
$ws.Range("H132").Value = 2259.8462
$ws.Range("I132").Value = 1480.871
$ws.Range("J132").Value = 3409.762
$ws.Range("K132").Value = 4442.613
$ws.Range("L132").Value = 10229.286
$ws.Range("M132").Value = -1912.613
$ws.Range("N132").Value = -15289.286
$ws.Range("H133").Value = 37975
$ws.Range("J133").Value = 37975
$ws.Range("L133").Value = 37975
$ws.Range("N133").Value = -48095

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 23333.334
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 23333.334
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 23333.334
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -23963.334
$ws.Range("H73").Value = 23333.334
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 23333.334
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 23333.334
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -25517.334
$ws.Range("H119").Value = 48665
$ws.Range("J119").Value = 48665
$ws.Range("L119").Value = 48665
$ws.Range("N119").Value = -58341
$ws.Range("H121").Value = 42509
$ws.Range("J121").Value = 42509
$ws.Range("L121").Value = 42509
$ws.Range("N121").Value = -46003
$ws.Range("H132").Value = 1354.2941
$ws.Range("I132").Value = 1096.9445
$ws.Range("J132").Value = 2346.9285
$ws.Range("K132").Value = 3290.8335
$ws.Range("L132").Value = 7040.7855
$ws.Range("M132").Value = -760.8335000000002
$ws.Range("N132").Value = -12100.7855
$ws.Range("H138").Value = 45437.5
$ws.Range("J138").Value = 45437.5
$ws.Range("L138").Value = 45437.5
$ws.Range("N138").Value = -55717.5
